$wb = $excel.ActiveWorkbook

# =====================================================================
# Step 1: Insert the new "2022-Q4" summary row into "总计" (first sheet)
# =====================================================================
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# The row that is now row 3 ("2022-Q2", formerly row 2) still carries the
# template formatting for a data row; copy it onto the freshly inserted
# row 2 before writing the new values so styles stay consistent.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122) # xlPasteFormats
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 10
$summary.Cells.Item(2,4).Value = 2.69

# The former row 8 (2020-Q4) is now row 9; its running index (col A) advances by one.
$summary.Cells.Item(9,1).Value = 7

# =====================================================================
# Step 2: Insert a brand-new "2022-Q4" worksheet, right after "总计" (i.e.
# immediately before the existing "2022-Q2" sheet) and fill it with the
# quarterly fund-holdings table.
# =====================================================================
$q2Before = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($q2Before)
$ws.Name = "2022-Q4"

# NOTE: adding a sheet shifts every subsequent sheets index, which can
# leave previously-grabbed COM references stale, so re-resolve "2022-Q2"
# by name now that the new sheet exists before using it as a copy source.
$q2 = $wb.Worksheets.Item("2022-Q2")

# Borrow the header/data formatting from "2022-Q2" so styles match the
# other quarterly sheets (bold/centered/bordered header, boxed index col).
$q2.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats
$q2.Range("A2:H2").Copy()
$ws.Range("A2:H11").PasteSpecial(-4122) # xlPasteFormats, tiled across all 10 data rows

$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Row 2: fund 516970
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "516970"
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "广发中证基建工程ETF"
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "73.53"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "99.57"
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = "2.42"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "1.7794"
$ws.Cells.Item(2,8).Value = 10

# Row 3: fund 006682
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,2).Value = "006682"
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "景顺长城中证500指数增强A"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "17.02"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "93.72"
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = "1.68"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "0.2859"
$ws.Cells.Item(3,8).Value = 10

# Row 4: fund 165525
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = "165525"
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "信诚中证基建工程指数（LOF）"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "9.51"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "94.24"
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = "2.29"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "0.2178"
$ws.Cells.Item(4,8).Value = 10

# Row 5: fund 000978
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,2).Value = "000978"
$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,3).Value = "景顺长城量化精选股票"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "7.44"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "94.11"
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = "2.05"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "0.1525"
$ws.Cells.Item(5,8).Value = 6

# Row 6: fund 001917
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).NumberFormat = "@"
$ws.Cells.Item(6,2).Value = "001917"
$ws.Cells.Item(6,3).NumberFormat = "@"
$ws.Cells.Item(6,3).Value = "招商量化精选股票A"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "5.91"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "94.08"
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = "1.44"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "0.0851"
$ws.Cells.Item(6,8).Value = 5

# Row 7: fund 007950
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).NumberFormat = "@"
$ws.Cells.Item(7,2).Value = "007950"
$ws.Cells.Item(7,3).NumberFormat = "@"
$ws.Cells.Item(7,3).Value = "招商量化精选股票C"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "5.28"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "94.08"
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = "1.44"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "0.0760"
$ws.Cells.Item(7,8).Value = 5

# Row 8: fund 013082
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).NumberFormat = "@"
$ws.Cells.Item(8,2).Value = "013082"
$ws.Cells.Item(8,3).NumberFormat = "@"
$ws.Cells.Item(8,3).Value = "信诚中证基建工程指数（LOF）C"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "2.86"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "94.24"
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = "2.29"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "0.0655"
$ws.Cells.Item(8,8).Value = 10

# Row 9: fund 009927
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).NumberFormat = "@"
$ws.Cells.Item(9,2).Value = "009927"
$ws.Cells.Item(9,3).NumberFormat = "@"
$ws.Cells.Item(9,3).Value = "工银瑞信聚利18个月定期开放混合A"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "2.18"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "26.90"
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = "1.21"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "0.0264"
$ws.Cells.Item(9,8).Value = 10

# Row 10: fund 009928
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).NumberFormat = "@"
$ws.Cells.Item(10,2).Value = "009928"
$ws.Cells.Item(10,3).NumberFormat = "@"
$ws.Cells.Item(10,3).Value = "工银瑞信聚利18个月定期开放混合C"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.39"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "26.90"
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = "1.21"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "0.0047"
$ws.Cells.Item(10,8).Value = 10

# Row 11: fund 016935
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).NumberFormat = "@"
$ws.Cells.Item(11,2).Value = "016935"
$ws.Cells.Item(11,3).NumberFormat = "@"
$ws.Cells.Item(11,3).Value = "景顺长城中证500指数增强C"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.00"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "93.72"
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = "1.68"
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 10

# Restore "General" display formatting on the text-typed columns without
# disturbing the values just written (the "@" format above was only needed
# transiently so Excel stored e.g. "73.53" as text instead of a number).
$ws.Range("B2:G11").Style = "Normal"

Write-Host "Added 2022-Q4 summary row and 2022-Q4 worksheet"
